# "Fixed import test cases - Task #967"
#
# Replaces the placeholder POW_IN/POW_SOMETHING/POW_OUT + BIII/BAAA/BUUU
# cost-equipment sample rows with real cost-type/name pairs, renames the
# sample battery equipment to BATTERY_1, and leaves the CostEquipments tab
# selected/active (instead of Header).
#
# NOTE: the cell-write order below matches the order the new shared
# strings first appear in the saved workbook (they are appended to the
# shared-string table in first-seen order), so it is kept exactly as is.

$wb = $excel.ActiveWorkbook

$wsHeader    = $wb.Worksheets.Item("Header")
$wsCostTypes = $wb.Worksheets.Item("CostTypes")
$wsCostEquip = $wb.Worksheets.Item("CostEquipments")

# --- CostEquipments sheet: update the cost-type / name pairs (rows 5-9) ---
$wsCostEquip.Range("C5").Value = "CostMaterial"
$wsCostEquip.Range("C6").Value = "CostPersonal"
$wsCostEquip.Range("C7").Value = "CostTest"
$wsCostEquip.Range("D7").Value = "Test"

# --- Header sheet: sample equipment name BATTERY -> BATTERY_1 ---
$wsHeader.Range("B6").Value = "BATTERY_1"

$wsCostEquip.Range("D6").Value = "Entwicklung"
$wsCostEquip.Range("D5").Value = "Material"

$wsCostEquip.Range("C8").Value = "CostTest"
$wsCostEquip.Range("D8").Value = "Test"

$wsCostEquip.Range("D9").Value = "Test"

# Row 8 gains a small custom row height
$wsCostEquip.Rows("8").RowHeight = 12.6

# --- Selections / active sheet ---
# CostTypes: selection moves from C5 to C16
$wsCostTypes.Range("C16").Select()

# CostEquipments: selection moves from C4 to B7
$wsCostEquip.Range("B7").Select()

# CostEquipments becomes the active / selected tab (was Header)
$wsCostEquip.Activate()
